$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.935.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.37%  '
$ws.Range("D3").Value = "'3.518.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.45%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'592.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'169.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.52%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'3.518.78"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.37%  '
$ws.Range("D9").Value = "'0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").Value = "'7.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("E11").Value = '  +5.57%  '
$ws.Range("D12").Value = "'0.438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").Value = "'4.125.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.46%  '
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").Value = "'28.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.24%  '
$ws.Range("D16").Value = "'66.860.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.15%  '
$ws.Range("E17").Value = '  +4.66%  '
$ws.Range("D18").Value = "'3.516.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.64%  '
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("D20").Value = "'14.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.10%  '
$ws.Range("D21").Value = "'390.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("D22").Value = "'7.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.74%  '
$ws.Range("D23").Value = "'73.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.11%  '
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("E25").Value = '  +2.93%  '
$ws.Range("E26").Value = '  +8.94%  '
$ws.Range("D27").Value = "'10.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.59%  '
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = "'6.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.55%  '
$ws.Range("E31").Value = '  +5.89%  '
$ws.Range("E32").Value = '  +4.51%  '
$ws.Range("D33").Value = "'23.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.12%  '
$ws.Range("D34").Value = "'7.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.18%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = "'1.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.79%  '
$ws.Range("D37").Value = "'161.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.75%  '
$ws.Range("D38").Value = "'0.910"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.32%  '
$ws.Range("D39").Value = "'1.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.94%  '
$ws.Range("D40").Value = "'0.0748"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.88%  '
$ws.Range("D41").Value = "'4.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.20%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'6.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.02%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'26.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("D44").Value = "'2.810.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = "'43.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.46%  '
$ws.Range("D46").Value = "'26.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.10%  '
$ws.Range("E47").Value = '  +10.53%  '
$ws.Range("D48").Value = "'0.0314"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.04%  '
$ws.Range("D49").Value = "'355.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.61%  '
$ws.Range("E50").Value = '  +6.70%  '
$ws.Range("D51").Value = "'33.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.29%  '
